# Implement cash flows for efficient building component rebate policy (#4)
#
# The "About" sheet cites the source paper for this variable via a
# hyperlinked URL in cell B6. The CEPE working paper has moved to a new
# URL on ethz.ch, so update the displayed/stored text of that cell to
# point at the new location.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

$ws.Range("B6").Value = "https://ethz.ch/content/dam/ethz/special-interest/mtec/cepe/cepe-dam/documents/research/cepe-wp/CEPE_WP86.pdf"
